# Add the personal website link right after every occurrence of the
# email address, on its own line (a manual line break followed by a
# new run carrying the same "muted" contact-info formatting), matching
# the contact-info pattern already used for the e-mail line.

$d = $word.ActiveDocument

$emailText = "BenniAustinDev@gmail.com"
$urlText   = "https://benniaustindev.github.io/"

$rng = $d.Content
$guard = 0

while ($rng.Find.Execute($emailText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) -and $guard -lt 50) {
    $guard = $guard + 1

    # Grab the formatting (rPr) of the matched e-mail run by copying its
    # FormattedText - this preserves color/sz/szCs exactly.
    $matchLen = $rng.End - $rng.Start
    $fmt = $rng.FormattedText

    # Move to just after the e-mail text (still before the pre-existing
    # trailing <w:br/> run) and insert a manual line break as its own run.
    $rng.Collapse(0)
    $insertPos = $rng.End
    $rng.InsertBreak(6)

    # Paste a copy of the e-mail run (same rPr) right after the new break,
    # then overwrite its text with the URL - this yields a distinct run
    # with matching color/sz/szCs instead of merging into a neighboring run.
    $dest = $d.Range($insertPos + 1, $insertPos + 1)
    $dest.FormattedText = $fmt

    $urlRng = $d.Range($insertPos + 1, $insertPos + 1 + $matchLen)
    $urlRng.Text = $urlText

    # Continue searching after the text we just inserted.
    $newEnd = $insertPos + 1 + $urlText.Length
    $rng = $d.Range($newEnd, $newEnd)
}
